$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.865.21'
$ws.Range('E2').Value = '  -3.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.907.93'
$ws.Range('E3').Value = '  -4.12%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.85'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.61'
$ws.Range('E6').Value = '  -5.98%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.907.49'
$ws.Range('E9').Value = '  -4.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.69'
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('E11').Value = '  -4.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.444'
$ws.Range('E12').Value = '  -3.95%  '
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.47'
$ws.Range('E14').Value = '  -6.18%  '
$ws.Range('E15').Value = '  +1.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.388.71'
$ws.Range('E16').Value = '  -4.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.809.82'
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.68'
$ws.Range('E18').Value = '  -5.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.907.55'
$ws.Range('E19').Value = '  -4.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '428.53'
$ws.Range('E20').Value = '  -5.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.52'
$ws.Range('E21').Value = '  -5.35%  '
$ws.Range('E22').Value = '  -2.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.10'
$ws.Range('E23').Value = '  -5.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.84'
$ws.Range('E24').Value = '  -2.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.92'
$ws.Range('E25').Value = '  -4.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.23'
$ws.Range('E26').Value = '  -4.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.93'
$ws.Range('E27').Value = '  -4.12%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.20'
$ws.Range('E30').Value = '  -3.80%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.23'
$ws.Range('E31').Value = '  -4.30%  '
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.43'
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.106'
$ws.Range('E34').Value = '  -3.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0866'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('E36').Value = '  -3.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.58'
$ws.Range('E37').Value = '  -5.82%  '
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.46'
$ws.Range('E39').Value = '  -1.97%  '
$ws.Range('E40').Value = '  -5.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.99'
$ws.Range('E41').Value = '  -5.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.59'
$ws.Range('E42').Value = '  -5.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.297'
$ws.Range('E43').Value = '  -4.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.30'
$ws.Range('E44').Value = '  -7.73%  '
$ws.Range('E45').Value = '  -2.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '375.59'
$ws.Range('E46').Value = '  -4.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.695.38'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.48'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.06'
$ws.Range('E50').Value = '  -5.79%  '
$ws.Range('E51').Value = '  -2.66%  '
